$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Update price (D) and volume (E) values for existing rows ---

Set-TextValue $ws.Range("D2") '26.861.01'
$ws.Range("E2").Value = '  +0.32%  '

Set-TextValue $ws.Range("D3") '1.642.20'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  -0.56%  '

Set-TextValue $ws.Range("D5") '218.33'
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("E6").Value = '  -0.67%  '

Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  -0.39%  '

Set-TextValue $ws.Range("D9") '0.0623'
$ws.Range("E9").Value = '  -1.08%  '

$ws.Range("E10").Value = '  +0.42%  '

Set-TextValue $ws.Range("D11") '0.0843'
$ws.Range("E11").Value = '  +0.08%  '

Set-TextValue $ws.Range("D12") '1.871.19'
$ws.Range("E12").Value = '  -0.05%  '

Set-TextValue $ws.Range("D13") '1.642.08'
$ws.Range("E13").Value = '  -0.55%  '

Set-TextValue $ws.Range("D14") '4.16'
$ws.Range("E14").Value = '  -0.19%  '

$ws.Range("E15").Value = '  +0.01%  '

Set-TextValue $ws.Range("D16") '65.28'
$ws.Range("E16").Value = '  +1.25%  '

Set-TextValue $ws.Range("D17") '26.847.78'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("E18").Value = '  -1.17%  '

Set-TextValue $ws.Range("D19") '216.71'
$ws.Range("E19").Value = '  +1.42%  '

Set-TextValue $ws.Range("D20") '1.00'
$ws.Range("E20").Value = '  -0.47%  '

Set-TextValue $ws.Range("D21") '4.37'
$ws.Range("E21").Value = '  -0.15%  '

Set-TextValue $ws.Range("D22") '6.58'
$ws.Range("E22").Value = '  +5.10%  '

Set-TextValue $ws.Range("D23") '2.37'
$ws.Range("E23").Value = '  -3.72%  '

$ws.Range("E24").Value = '  -1.62%  '

Set-TextValue $ws.Range("D25") '147.52'
$ws.Range("E25").Value = '  +1.66%  '

Set-TextValue $ws.Range("D26") '1.01'
$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("E28").Value = '  +1.04%  '

Set-TextValue $ws.Range("D29") '15.75'
$ws.Range("E29").Value = '  +0.63%  '

Set-TextValue $ws.Range("D30") '0.0508'
$ws.Range("E30").Value = '  -0.44%  '

Set-TextValue $ws.Range("D31") '1.20'
$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("E33").Value = '  +0.15%  '

Set-TextValue $ws.Range("D34") '1.281.29'
$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("E35").Value = '  +0.84%  '

$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("E37").Value = '  -0.84%  '

Set-TextValue $ws.Range("D38") '0.532'
$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("E39").Value = '  -0.50%  '

Set-TextValue $ws.Range("D40") '1.00'
$ws.Range("E40").Value = '  -0.33%  '

Set-TextValue $ws.Range("D41") '0.804'
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("E42").Value = '  -0.27%  '

Set-TextValue $ws.Range("D43") '1.782.18'
$ws.Range("E43").Value = '  -0.68%  '

$ws.Range("E44").Value = '  -6.16%  '

Set-TextValue $ws.Range("D45") '92.65'
$ws.Range("E45").Value = '  +1.17%  '

Set-TextValue $ws.Range("D46") '61.19'
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("E47").Value = '  -1.14%  '

# --- New coin (BabyDogeCoin) inserted at row 48, shifting Cronos/EnergySwap/
# --- Algorand down a row; USDD (previously row 51) drops off the bottom of the list ---

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D48") '0.0₆0102'
$ws.Range("E48").Value = '  -2.56%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D49") '0.0516'
$ws.Range("E49").Value = '  -1.72%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D50") '7.58'
$ws.Range("E50").Value = '  -1.70%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.0966'
$ws.Range("E51").Value = '  -0.94%  '

